# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables pseudo-attributes embedded in the header cell text of each
# sheet are renamed from UpperCamelCase to lowerCamelCase:
#   ObjTablesVersion -> objTablesVersion
#   Type             -> type
#   Id               -> id

$wb = $excel.ActiveWorkbook

$wsNormal = $wb.Worksheets.Item("!!Normals")
$wsNormal.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsNormal.Range("A2").Value = "!!ObjTables type='Data' id='Normal'"

$wsTransposed = $wb.Worksheets.Item("!!Transposed")
$wsTransposed.Range("A1").Value = "!!ObjTables type='Data' id='Transposed'"
